$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.372.17'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.593.41'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.507'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.49'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '1.816.47'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.602.70'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.525'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.78'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '26.363.79'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").Value = '0.0₃0733'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '212.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.44%  '
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.21%  '
$ws.Range("E24").Value = '  -2.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").Value = '1.337.03'
$ws.Range("E34").Value = '  +4.05%  '
$ws.Range("E35").Value = '  -1.26%  '
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.820'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.81%  '
$ws.Range("E41").Value = '  -0.39%  '
$ws.Range("E42").Value = '  -23.38%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.767'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '1.729.81'
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("E48").Value = '  -3.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0985'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0505'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.57%  '
